# Update crypto price/volume data (and reorder the last three coin rows)
# as produced by the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.042.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.240.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.34%  '

$ws.Range("E7").Value = '  -2.78%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -6.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0823'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.47%  '

$ws.Range("E13").Value = '  -3.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.581.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.843'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.245.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.889.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0982'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.07%  '

$ws.Range("E21").Value = '  -3.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.27%  '

$ws.Range("E24").Value = '  -7.04%  '

$ws.Range("E25").Value = '  -8.05%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.82%  '

$ws.Range("E29").Value = '  -4.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.56'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0836'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.94%  '

$ws.Range("E35").Value = '  -3.89%  '

$ws.Range("E36").Value = '  -4.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.89'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.39%  '

$ws.Range("E38").Value = '  -2.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.24'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -12.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0309'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.05%  '

$ws.Range("E43").Value = '  +0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.709.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '83.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.71%  '

$ws.Range("E46").Value = '  -6.44%  '

$ws.Range("E47").Value = '  -5.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.42%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '56.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.61%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.30%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.61'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.89%  '
